$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '70.503.94'
$ws.Range('E2').Value = '  -3.56%  '

# Row 3
$ws.Range('D3').Value = '3.840.29'
$ws.Range('E3').Value = '  -4.07%  '

# Row 4
$ws.Range('E4').Value = '  -0.04%  '

# Row 5
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '593.25'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +0.16%  '

# Row 6
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '165.98'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +3.28%  '

# Row 7
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.670'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  -2.72%  '

# Row 8
$ws.Range('E8').Value = '  +0.21%  '

# Row 9
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.745'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  -1.06%  '

# Row 10
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.174'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  +3.24%  '

# Row 11
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '52.92'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -2.37%  '

# Row 12
$ws.Range('E12').Value = '  -0.37%  '

# Row 13
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '11.32'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +2.80%  '

# Row 14
$ws.Range('D14').Value = '4.452.46'
$ws.Range('E14').Value = '  -3.73%  '

# Row 15
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '20.97'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  +2.53%  '

# Row 16
$ws.Range('D16').Value = '3.851.22'
$ws.Range('E16').Value = '  -3.52%  '

# Row 17
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '13.78'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  -2.72%  '

# Row 18
$ws.Range('E18').Value = '  -6.01%  '

# Row 19
$ws.Range('E19').Value = '  -2.18%  '

# Row 20
$ws.Range('D20').Value = '70.329.92'
$ws.Range('E20').Value = '  -3.41%  '

# Row 21
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '435.06'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -0.47%  '

# Row 22
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '4.70'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -1.68%  '

# Row 23
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '93.66'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -2.87%  '

# Row 24
$ws.Range('E24').Value = '  -5.81%  '

# Row 25
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '13.84'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -3.43%  '

# Row 26
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '11.16'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -1.44%  '

# Row 27
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '3.95'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -12.02%  '

# Row 28
$ws.Range('E28').Value = '  -0.05%  '

# Row 29
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '10.41'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  -1.14%  '

# Row 30
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '35.01'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -4.04%  '

# Row 31
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '8.10'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +1.81%  '

# Row 32
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '13.42'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -2.28%  '

# Row 33
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '47.96'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -1.71%  '

# Row 34
$ws.Range('E34').Value = '  -5.26%  '

# Row 35
$ws.Range('E35').Value = '  -0.69%  '

# Row 36
$ws.Range('D36').Value = '0.0₃0984'
$ws.Range('E36').Value = '  +11.24%  '

# Row 37
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '635.64'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -5.62%  '

# Row 38
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.422'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -3.68%  '

# Row 39
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +0.04%  '

# Row 40
$ws.Range('E40').Value = '  -1.25%  '

# Row 41
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +0.09%  '

# Row 42
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '3.26'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -3.39%  '

# Row 43
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '3.21'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +20.83%  '

# Row 44
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '0.0467'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -4.54%  '

# Row 45
$ws.Range('B45').Value = 'Fetch.AI'
$ws.Range('C45').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '2.72'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +3.15%  '

# Row 46
$ws.Range('B46').Value = 'THORChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '9.98'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -8.36%  '

# Row 47
$ws.Range('E47').Value = '  -4.62%  '

# Row 48
$ws.Range('E48').Value = '  -15.13%  '

# Row 49
$ws.Range('D49').Value = '2.846.33'
$ws.Range('E49').Value = '  +0.71%  '

# Row 50
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '3.25'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -4.73%  '

# Row 51
$ws.Range('E51').Value = '  +0.19%  '

